$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "98.533.92"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.340.36"
$ws.Range("E3").Value = "  -0.93%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "257.74"
$ws.Range("E5").Value = "  -0.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "663.75"
$ws.Range("E6").Value = "  +5.58%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.54"
$ws.Range("E7").Value = "  +12.16%  "
$ws.Range("E8").Value = "  +16.79%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.09"
$ws.Range("E9").Value = "  +27.02%  "
$ws.Range("E10").Value = "  -0.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.336.40"
$ws.Range("E11").Value = "  -0.99%  "
$ws.Range("E12").Value = "  +5.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "42.92"
$ws.Range("E13").Value = "  +19.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000266"
$ws.Range("E14").Value = "  +7.86%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "98.220.81"
$ws.Range("E15").Value = "  -0.86%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.963.31"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("E17").Value = "  +2.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.336.93"
$ws.Range("E18").Value = "  -0.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.52"
$ws.Range("E19").Value = "  +23.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.90"
$ws.Range("E20").Value = "  +10.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.61"
$ws.Range("E21").Value = "  +1.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "526.01"
$ws.Range("E22").Value = "  +6.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.43"
$ws.Range("E23").Value = "  +11.77%  "
$ws.Range("B24").Value = "Stellar"
$ws.Range("C24").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.439"
$ws.Range("E24").Value = "  +57.67%  "
$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000210"
$ws.Range("E25").Value = "  +0.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.33"
$ws.Range("E26").Value = "  +12.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "101.48"
$ws.Range("E27").Value = "  +14.78%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.69"
$ws.Range("E28").Value = "  +6.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.519.40"
$ws.Range("E29").Value = "  -0.83%  "
$ws.Range("E30").Value = "  +17.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.05"
$ws.Range("E32").Value = "  +16.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.192"
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "29.25"
$ws.Range("E35").Value = "  +4.41%  "
$ws.Range("E36").Value = "  +17.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.88"
$ws.Range("E37").Value = "  +7.76%  "
$ws.Range("B38").Value = "PancakeSwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.12"
$ws.Range("E38").Value = "  +8.19%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.159"
$ws.Range("E39").Value = "  +5.63%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "532.13"
$ws.Range("E40").Value = "  +6.83%  "
$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.34"
$ws.Range("E41").Value = "  +6.36%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0450"
$ws.Range("E42").Value = "  +37.71%  "
$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "24.70"
$ws.Range("E43").Value = "  -0.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.76"
$ws.Range("E44").Value = "  -3.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.821"
$ws.Range("E45").Value = "  +5.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.35"
$ws.Range("E46").Value = "  +2.83%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("E48").Value = "  +6.88%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.91"
$ws.Range("E49").Value = "  +21.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.13"
$ws.Range("E50").Value = "  +10.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "164.83"
$ws.Range("E51").Value = "  +3.02%  "
